# Pequena alteracao sobre os testes
#
# The paragraph that begins "Os testes em cada fase do ciclo de vida..."
# is a single run. The edit turns it into six runs (same font settings
# throughout) by:
#   1. turning " e sem pressa" into ", sem pressa"
#   2. inserting " e testando todas as possibilidades de entrada de
#      dados" right after "disponível)" (and before the existing ". Por
#      mais")
#   3. inserting ", dando um passo de cada vez" right before the final
#      period of the paragraph
#
# iron_native's Range.Text setter rebuilds/merges runs with identical
# rPr, so every edit is immediately followed by a harmless Bold
# true->false round trip on the same Range. That forces the engine to
# keep a hard run boundary at that Range even though the formatting
# ends up unchanged - exactly mirroring the run layout in the target
# XML. All the Range.Text assignments are done first and only then are
# the boundaries "locked in", otherwise a later Text assignment
# re-merges earlier, already-split runs.

function Lock-RunBoundary($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

$d = $word.ActiveDocument

# Locate the target paragraph robustly (rather than hard-coding an index).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i).Range
    if ($candidate.Text -like "Os testes em cada fase do ciclo de vida*") {
        $target = $candidate
        break
    }
}

$para = $target.Duplicate

# --- Edit 1: " e" -> "," right before " sem pressa" -----------------
$findRange = $para.Duplicate
$findRange.Find.Execute(" e sem pressa")
$commaRange = $findRange.Duplicate
$commaRange.MoveEnd(1, -11)                 # keep just the " e" (2 chars)
$commaRange.Text = ","

# --- Edit 2: insert the new sentence about testing inputs ------------
$findRange2 = $para.Duplicate
$findRange2.Find.Execute("disponível). Por mais")
$insertRange = $findRange2.Duplicate
$insertRange.MoveStart(1, 11)               # skip over "disponível)"
$insertRange.MoveEnd(1, -9)                 # drop the trailing " Por mais", keep just "."
$insertRange.Text = " e testando todas as possibilidades de entrada de dados."

# --- Edit 3: insert the closing remark before the final period -------
$findRange3 = $para.Duplicate
$findRange3.Find.Execute("como um todo.")
$finalRange = $findRange3.Duplicate
$finalRange.MoveStart(1, 12)                # skip over "como um todo", keep just "."
$finalRange.Text = ", dando um passo de cada vez."

# --- Now lock in the run boundaries for all three edits ---------------
Lock-RunBoundary $commaRange
Lock-RunBoundary $insertRange
Lock-RunBoundary $finalRange
